$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.614.14'
$ws.Range('E2').Value = '  +1.43%  '
$ws.Range('D3').Value = '3.323.20'
$ws.Range('E3').Value = '  +2.41%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '580.40'
$ws.Range('E5').Value = '  +0.31%  '
$ws.Range('D6').Value = '174.41'
$ws.Range('E6').Value = '  +1.07%  '
$ws.Range('E7').Value = '  +0.13%  '
$ws.Range('E8').Value = '  +2.35%  '
$ws.Range('D9').Value = '3.319.09'
$ws.Range('E9').Value = '  +2.55%  '
$ws.Range('D10').Value = '0.181'
$ws.Range('E10').Value = '  +6.09%  '
$ws.Range('D11').Value = '0.578'
$ws.Range('E11').Value = '  +1.85%  '
$ws.Range('D12').Value = '46.79'
$ws.Range('E12').Value = '  +4.95%  '
$ws.Range('E13').Value = '  +1.26%  '
$ws.Range('D14').Value = '691.42'
$ws.Range('E14').Value = '  +5.15%  '
$ws.Range('D15').Value = '3.868.70'
$ws.Range('E15').Value = '  +2.66%  '
$ws.Range('D16').Value = '8.36'
$ws.Range('E16').Value = '  +1.81%  '
$ws.Range('D17').Value = '67.628.66'
$ws.Range('E17').Value = '  +1.49%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = '0.119'
$ws.Range('E18').Value = '  +0.53%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.327.71'
$ws.Range('E19').Value = '  +2.54%  '
$ws.Range('D20').Value = '17.48'
$ws.Range('E20').Value = '  +1.87%  '
$ws.Range('D21').Value = '11.05'
$ws.Range('E21').Value = '  +3.55%  '
$ws.Range('D22').Value = '0.889'
$ws.Range('E22').Value = '  +1.83%  '
$ws.Range('D23').Value = '5.49'
$ws.Range('E23').Value = '  +3.82%  '
$ws.Range('D24').Value = '16.81'
$ws.Range('E24').Value = '  -0.48%  '
$ws.Range('D25').Value = '101.09'
$ws.Range('E25').Value = '  +4.38%  '
$ws.Range('D26').Value = '3.89'
$ws.Range('E26').Value = '  +1.61%  '
$ws.Range('E27').Value = '  +1.99%  '
$ws.Range('E28').Value = '  +3.53%  '
$ws.Range('D29').Value = '32.80'
$ws.Range('E29').Value = '  +2.85%  '
$ws.Range('D30').Value = '8.49'
$ws.Range('E30').Value = '  +2.70%  '
$ws.Range('E31').Value = '  +2.97%  '
$ws.Range('D32').Value = '569.80'
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('D33').Value = '10.96'
$ws.Range('E33').Value = '  +1.15%  '
$ws.Range('E34').Value = '  +3.23%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('B36').Value = 'OKB'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D36').Value = '57.24'
$ws.Range('E36').Value = '  +3.48%  '
$ws.Range('B37').Value = 'Maker'
$ws.Range('C37').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D37').Value = '3.714.11'
$ws.Range('E37').Value = '  -1.09%  '
$ws.Range('D38').Value = '3.26'
$ws.Range('E38').Value = '  -3.55%  '
$ws.Range('D39').Value = '35.00'
$ws.Range('E39').Value = '  +9.34%  '
$ws.Range('E40').Value = '  +3.56%  '
$ws.Range('D41').Value = '3.14'
$ws.Range('E41').Value = '  +5.37%  '
$ws.Range('D42').Value = '2.60'
$ws.Range('E42').Value = '  +0.48%  '
$ws.Range('E43').Value = '  +3.34%  '
$ws.Range('D44').Value = '0.333'
$ws.Range('E44').Value = '  +3.15%  '
$ws.Range('D45').Value = '0.0₃0665'
$ws.Range('E45').Value = '  +2.10%  '
$ws.Range('E46').Value = '  +2.14%  '
$ws.Range('D47').Value = '2.63'
$ws.Range('E47').Value = '  +3.98%  '
$ws.Range('E48').Value = '  +1.72%  '
$ws.Range('E49').Value = '  +0.09%  '
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('D51').Value = '131.12'
$ws.Range('E51').Value = '  +2.32%  '
